$wb = $excel.ActiveWorkbook

# Sheet ALC row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 309.16666
$ws.Range("I55").Value = 181.11111
$ws.Range("J55").Value = 693.3333
$ws.Range("K55").Value = 181.11111
$ws.Range("L55").Value = 693.3333
$ws.Range("M55").Value = 32.88889
$ws.Range("N55").Value = -1121.3333

# Sheet ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6421.759
$ws.Range("I74").Value = 5338.077
$ws.Range("K74").Value = 5338.077
$ws.Range("M74").Value = -4402.077

# Sheet ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 6421.759
$ws.Range("I77").Value = 5338.077
$ws.Range("K77").Value = 26690.385
$ws.Range("M77").Value = -22010.385

# Sheet ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3380.8948
$ws.Range("J98").Value = 5745
$ws.Range("L98").Value = 5745
$ws.Range("N98").Value = -8741

# Sheet ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 3682.6924
$ws.Range("I101").Value = 1684.875
$ws.Range("K101").Value = 5054.625
$ws.Range("M101").Value = -3432.625

# Sheet ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3380.8948
$ws.Range("J122").Value = 5745
$ws.Range("L122").Value = 17235
$ws.Range("N122").Value = -22135

# Sheet ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2094.5
$ws.Range("I63").Value = 2094.5
$ws.Range("K63").Value = 2094.5
$ws.Range("M63").Value = -1408.5

# Sheet ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2094.5
$ws.Range("I66").Value = 2094.5
$ws.Range("K66").Value = 10472.5
$ws.Range("M66").Value = -7040.5

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3090.2222
$ws.Range("I74").Value = 1479.1428
$ws.Range("K74").Value = 1479.1428
$ws.Range("M74").Value = -605.1428000000001

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3090.2222
$ws.Range("I77").Value = 1479.1428
$ws.Range("K77").Value = 7395.714
$ws.Range("M77").Value = -3027.714

# Sheet ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2079.8333
$ws.Range("I102").Value = 1721.75
$ws.Range("K102").Value = 1721.75
$ws.Range("M102").Value = -99.75

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2675.1365
$ws.Range("I132").Value = 2136.8
$ws.Range("K132").Value = 6410.400000000001
$ws.Range("M132").Value = -3880.400000000001

# Sheet BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3414.7778
$ws.Range("I99").Value = 3504.7144
$ws.Range("J99").Value = 3100
$ws.Range("K99").Value = 3504.7144
$ws.Range("L99").Value = 3100
$ws.Range("M99").Value = -2006.7144
$ws.Range("N99").Value = -6096

# Sheet CRP row 81
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()

# Sheet CRP row 82
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Sheet CRP row 84
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()

# Sheet CRP row 85
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6592.273
$ws.Range("J134").Value = 7699.4
$ws.Range("L134").Value = 23098.2
$ws.Range("N134").Value = -28168.2

# Sheet CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 7500.8335
$ws.Range("J58").Value = 10335
$ws.Range("L58").Value = 31005
$ws.Range("N58").Value = -31261

# Sheet CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 3342.25
$ws.Range("I69").Value = 3276.8572
$ws.Range("K69").Value = 9830.571599999999
$ws.Range("M69").Value = -9019.571599999999

# Sheet CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 3342.25
$ws.Range("I72").Value = 3276.8572
$ws.Range("K72").Value = 29491.7148
$ws.Range("M72").Value = -25435.7148

# Sheet GSM row 38
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 27497.75
$ws.Range("J38").Value = 27497.75
$ws.Range("L38").Value = 27497.75
$ws.Range("N38").Value = -28423.75

# Sheet GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 32180
$ws.Range("J57").Value = 39600
$ws.Range("L57").Value = 39600
$ws.Range("N57").Value = -41240

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1499.6666
$ws.Range("I80").Value = 1249.5
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1249.5
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -251.5
$ws.Range("N80").Value = -3996

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1499.6666
$ws.Range("I83").Value = 1249.5
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 6247.5
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -1255.5
$ws.Range("N83").Value = -19984

# Sheet GSM row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Sheet LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2361.875
$ws.Range("I93").Value = 2361.875
$ws.Range("K93").Value = 2361.875
$ws.Range("M93").Value = -1113.875

# Sheet WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20499.666
$ws.Range("I4").Value = 24199.6
$ws.Range("K4").Value = 24199.6
$ws.Range("M4").Value = -24086.6

# Sheet WVR row 6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1671.3334
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 10003
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 10003
$ws.Range("M6").Value = 110
$ws.Range("N6").Value = -10233

# Sheet WVR row 9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 9999
$ws.Range("I9").Value = 9999
$ws.Range("K9").Value = 9999
$ws.Range("M9").Value = -9859

# Sheet WVR row 52
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 11395.429
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7927.25
$ws.Range("I81").Value = 7927.25
$ws.Range("K81").Value = 15854.5
$ws.Range("M81").Value = -14793.5

# Sheet WVR row 82
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 110000
$ws.Range("J82").Value = 110000
$ws.Range("L82").Value = 110000
$ws.Range("N82").Value = -110766

# Sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7927.25
$ws.Range("I84").Value = 7927.25
$ws.Range("K84").Value = 79272.5
$ws.Range("M84").Value = -73968.5

# Sheet WVR row 85
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 110000
$ws.Range("J85").Value = 110000
$ws.Range("L85").Value = 110000
$ws.Range("N85").Value = -112652

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2877.0667
$ws.Range("I132").Value = 2759.182
$ws.Range("K132").Value = 8277.545999999998
$ws.Range("M132").Value = -5747.545999999998
